$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'330.59"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'-0.03%"
$ws.Range("E2").Style = "Normal"
$ws.Range("G2").Value = "'20"
$ws.Range("G2").Style = "Normal"
$ws.Range("D3").Value = "'41.58"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'0.93%"
$ws.Range("E3").Style = "Normal"
$ws.Range("G3").Value = "'20"
$ws.Range("G3").Style = "Normal"
$ws.Range("D4").Value = "'5.706"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'0.34%"
$ws.Range("E4").Style = "Normal"
$ws.Range("G4").Value = "'20"
$ws.Range("G4").Style = "Normal"
$ws.Range("D5").Value = "'0.08393"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'4.24%"
$ws.Range("E5").Style = "Normal"
$ws.Range("G5").Value = "'20"
$ws.Range("G5").Style = "Normal"
$ws.Range("D6").Value = "'8.820"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'1.06%"
$ws.Range("E6").Style = "Normal"
$ws.Range("G6").Value = "'20"
$ws.Range("G6").Style = "Normal"
$ws.Range("D7").Value = "'2.002"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'-1.03%"
$ws.Range("E7").Style = "Normal"
$ws.Range("G7").Value = "'20"
$ws.Range("G7").Style = "Normal"
$ws.Range("D8").Value = "'4.486"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'-1.11%"
$ws.Range("E8").Style = "Normal"
$ws.Range("G8").Value = "'20"
$ws.Range("G8").Style = "Normal"
$ws.Range("G9").Value = "'20"
$ws.Range("G9").Style = "Normal"
$ws.Range("D10").Value = "'0.9238"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'0.40%"
$ws.Range("E10").Style = "Normal"
$ws.Range("G10").Value = "'20"
$ws.Range("G10").Style = "Normal"
$ws.Range("D11").Value = "'0.1285"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'2.39%"
$ws.Range("E11").Style = "Normal"
$ws.Range("G11").Value = "'20"
$ws.Range("G11").Style = "Normal"
$ws.Range("E12").Value = "'1.68%"
$ws.Range("E12").Style = "Normal"
$ws.Range("G12").Value = "'20"
$ws.Range("G12").Style = "Normal"
$ws.Range("D13").Value = "'0.09466"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'1.47%"
$ws.Range("E13").Style = "Normal"
$ws.Range("G13").Value = "'20"
$ws.Range("G13").Style = "Normal"
$ws.Range("E14").Value = "'3.15%"
$ws.Range("E14").Style = "Normal"
$ws.Range("G14").Value = "'20"
$ws.Range("G14").Style = "Normal"
$ws.Range("D15").Value = "'0.1062"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'0.96%"
$ws.Range("E15").Style = "Normal"
$ws.Range("G15").Value = "'20"
$ws.Range("G15").Style = "Normal"
$ws.Range("D16").Value = "'0.001312"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'0.73%"
$ws.Range("E16").Style = "Normal"
$ws.Range("G16").Value = "'20"
$ws.Range("G16").Style = "Normal"
$ws.Range("D17").Value = "'0.006111"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'-2.53%"
$ws.Range("E17").Style = "Normal"
$ws.Range("G17").Value = "'20"
$ws.Range("G17").Style = "Normal"
$ws.Range("D18").Value = "'3.423"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'1.71%"
$ws.Range("E18").Style = "Normal"
$ws.Range("G18").Value = "'20"
$ws.Range("G18").Style = "Normal"
$ws.Range("E19").Value = "'0.71%"
$ws.Range("E19").Style = "Normal"
$ws.Range("G19").Value = "'20"
$ws.Range("G19").Style = "Normal"
$ws.Range("D20").Value = "'8.967"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'2.65%"
$ws.Range("E20").Style = "Normal"
$ws.Range("G20").Value = "'20"
$ws.Range("G20").Style = "Normal"
$ws.Range("E21").Value = "'-4.04%"
$ws.Range("E21").Style = "Normal"
$ws.Range("G21").Value = "'20"
$ws.Range("G21").Style = "Normal"
$ws.Range("D22").Value = "'0.2511"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'-5.62%"
$ws.Range("E22").Style = "Normal"
$ws.Range("G22").Value = "'20"
$ws.Range("G22").Style = "Normal"
$ws.Range("D23").Value = "'0.04405"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'-0.80%"
$ws.Range("E23").Style = "Normal"
$ws.Range("G23").Value = "'20"
$ws.Range("G23").Style = "Normal"
$ws.Range("D24").Value = "'0.001246"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'-1.49%"
$ws.Range("E24").Style = "Normal"
$ws.Range("G24").Value = "'20"
$ws.Range("G24").Style = "Normal"
$ws.Range("D25").Value = "'0.004363"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'1.92%"
$ws.Range("E25").Style = "Normal"
$ws.Range("G25").Value = "'20"
$ws.Range("G25").Style = "Normal"
$ws.Range("D26").Value = "'0.0001192"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'-4.29%"
$ws.Range("E26").Style = "Normal"
$ws.Range("G26").Value = "'20"
$ws.Range("G26").Style = "Normal"
$ws.Range("D27").Value = "'0.0003994"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'0.02%"
$ws.Range("E27").Style = "Normal"
$ws.Range("G27").Value = "'20"
$ws.Range("G27").Style = "Normal"
$ws.Range("G28").Value = "'20"
$ws.Range("G28").Style = "Normal"
$ws.Range("G29").Value = "'20"
$ws.Range("G29").Style = "Normal"
$ws.Range("G30").Value = "'20"
$ws.Range("G30").Style = "Normal"
$ws.Range("G31").Value = "'20"
$ws.Range("G31").Style = "Normal"
$ws.Range("G32").Value = "'20"
$ws.Range("G32").Style = "Normal"
$ws.Range("G33").Value = "'20"
$ws.Range("G33").Style = "Normal"
$ws.Range("G34").Value = "'20"
$ws.Range("G34").Style = "Normal"
$ws.Range("G35").Value = "'20"
$ws.Range("G35").Style = "Normal"
$ws.Range("G36").Value = "'20"
$ws.Range("G36").Style = "Normal"
$ws.Range("G37").Value = "'20"
$ws.Range("G37").Style = "Normal"
$ws.Range("G38").Value = "'20"
$ws.Range("G38").Style = "Normal"
$ws.Range("D39").Value = "'0.02839"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'-0.55%"
$ws.Range("E39").Style = "Normal"
$ws.Range("G39").Value = "'20"
$ws.Range("G39").Style = "Normal"
$ws.Range("E40").Value = "'0.96%"
$ws.Range("E40").Style = "Normal"
$ws.Range("G40").Value = "'20"
$ws.Range("G40").Style = "Normal"
$ws.Range("D41").Value = "'0.007953"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'2.21%"
$ws.Range("E41").Style = "Normal"
$ws.Range("G41").Value = "'20"
$ws.Range("G41").Style = "Normal"
$ws.Range("D42").Value = "'0.1434"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'1.68%"
$ws.Range("E42").Style = "Normal"
$ws.Range("G42").Value = "'20"
$ws.Range("G42").Style = "Normal"
$ws.Range("D43").Value = "'0.008982"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'-9.88%"
$ws.Range("E43").Style = "Normal"
$ws.Range("G43").Value = "'20"
$ws.Range("G43").Style = "Normal"
$ws.Range("D44").Value = "'0.002093"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'-1.68%"
$ws.Range("E44").Style = "Normal"
$ws.Range("G44").Value = "'20"
$ws.Range("G44").Style = "Normal"
$ws.Range("D45").Value = "'0.01173"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'-1.46%"
$ws.Range("E45").Style = "Normal"
$ws.Range("G45").Value = "'20"
$ws.Range("G45").Style = "Normal"
$ws.Range("D46").Value = "'0.00006937"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'2.22%"
$ws.Range("E46").Style = "Normal"
$ws.Range("G46").Value = "'20"
$ws.Range("G46").Style = "Normal"
$ws.Range("E47").Value = "'-0.28%"
$ws.Range("E47").Style = "Normal"
$ws.Range("G47").Value = "'20"
$ws.Range("G47").Style = "Normal"
$ws.Range("D48").Value = "'0.003464"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'14.55%"
$ws.Range("E48").Style = "Normal"
$ws.Range("G48").Value = "'20"
$ws.Range("G48").Style = "Normal"
$ws.Range("D49").Value = "'0.002280"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'-0.34%"
$ws.Range("E49").Style = "Normal"
$ws.Range("G49").Value = "'20"
$ws.Range("G49").Style = "Normal"
$ws.Range("D50").Value = "'0.00002103"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'-0.28%"
$ws.Range("E50").Style = "Normal"
$ws.Range("G50").Value = "'20"
$ws.Range("G50").Style = "Normal"
$ws.Range("D51").Value = "'0.0002003"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'-0.28%"
$ws.Range("E51").Style = "Normal"
$ws.Range("G51").Value = "'20"
$ws.Range("G51").Style = "Normal"
